# Apply the update described by the commit: the "casos confirmados" table
# (rows 2-99) and the "obitos" table (rows 100-132) are unified into a single
# continuous list. This means:
#   - the "cidades"/"Casos confirmados" sub-header row (row 2) is removed
#   - the "outros estados" (row 98) and "outros paises" (row 99) summary
#     rows are removed
#   - the "cidade"/"Obtos" sub-header row (row 100) is removed
# Every row below each deleted row shifts up to close the gap, yielding a
# final used range of A1:B128 (down from A1:B132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from bottom to top so row numbers of not-yet-deleted rows stay valid.
$ws.Rows.Item(100).Delete()
$ws.Rows.Item(99).Delete()
$ws.Rows.Item(98).Delete()
$ws.Rows.Item(2).Delete()
